# Applies weekly fruit/vegetable price update: shuffles Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M) and
# Precio $/Kg (P) values across rows 2-11 (row 7 left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedio, PrecioKg)
$data = @{
    2  = @(44203, 30, 2000, 2000, 2000, 2000)
    3  = @(44474, 20, 1600, 1600, 1600, 1600)
    4  = @(44447, 75, 2200, 2200, 2200, 2200)
    5  = @(44497, 50, 2200, 2200, 2200, 2200)
    6  = @(44484, 40, 2200, 2200, 2200, 2200)
    8  = @(44453, 20, 2300, 2300, 2300, 2300)
    9  = @(44483, 50, 2200, 2200, 2200, 2200)
    10 = @(44476, 30, 2200, 2200, 2200, 2200)
    11 = @(44496, 40, 2200, 2200, 2200, 2200)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]

    $ws.Cells.Item($row, 4).Value = $values[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $values[1]  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $values[2]  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $values[3]  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $values[4]  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $values[5]  # P - Precio $/Kg
}
